$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Update header date
$d.Content.Find.Execute("2025-08-19 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-08-20 Wednesday", 2) | Out-Null

# Row 1 (table row 1)
$t.Cell(1,1).Range.Text = "57÷4=14, 1"
$t.Cell(1,2).Range.Text = "79÷4=19, 3"
$t.Cell(1,3).Range.Text = "12÷6=2, 0"
$t.Cell(1,4).Range.Text = "59÷4=14, 3"
$t.Cell(1,5).Range.Text = "41÷9=4, 5"

# Row 5 (table row 5)
$t.Cell(5,1).Range.Text = "53÷6=8, 5"
$t.Cell(5,2).Range.Text = "20÷7=2, 6"
$t.Cell(5,3).Range.Text = "68÷3=22, 2"
$t.Cell(5,4).Range.Text = "37÷2=18, 1"
$t.Cell(5,5).Range.Text = "50÷7=7, 1"

# Row 9 (table row 9)
$t.Cell(9,1).Range.Text = "31÷6=5, 1"
$t.Cell(9,2).Range.Text = "96÷7=13, 5"
$t.Cell(9,3).Range.Text = "15÷4=3, 3"
$t.Cell(9,4).Range.Text = "48÷3=16, 0"
$t.Cell(9,5).Range.Text = "74÷9=8, 2"

# Row 13 (table row 13)
$t.Cell(13,1).Range.Text = "76÷9=8, 4"
$t.Cell(13,2).Range.Text = "81÷5=16, 1"
$t.Cell(13,3).Range.Text = "50÷5=10, 0"
$t.Cell(13,4).Range.Text = "71÷3=23, 2"
$t.Cell(13,5).Range.Text = "89÷2=44, 1"

# Row 17 (table row 17)
$t.Cell(17,1).Range.Text = "50÷5=10, 0"
$t.Cell(17,2).Range.Text = "23÷4=5, 3"
$t.Cell(17,3).Range.Text = "76÷3=25, 1"
$t.Cell(17,4).Range.Text = "25÷9=2, 7"
$t.Cell(17,5).Range.Text = "84÷4=21, 0"
